$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (the last data row) is being replaced from a specific named user
# to a generic "Guest" placeholder row.
$ws.Range("A12").Value = "U0000"
$ws.Range("B12").Value = "کاربر مهمان"
$ws.Range("C12").Value = "Guest"
$ws.Range("E12").Value = "Guest"

# Update the active selection to the whole of row 12.
$ws.Range("A12:E12").Select()
